$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vitreous concentrations")
$ws.Range("D76").Value = 35.5
Write-Host "done"
